$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates per diff. D-column values are forced to Text
# (NumberFormat '@') before assignment and the style reset to
# 'Normal' afterwards, because several price strings look like
# plain numbers (e.g. '215.80', '19.64') and Excel's COM Value
# setter auto-converts those to doubles otherwise -- losing
# trailing zeros / exact text formatting that the source sheet
# relies on (prices are stored as inline/shared strings, not numbers).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.053.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.645.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.507"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("E8").Value = "  +0.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.256"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.64"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0797"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.47%  "
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.622.40"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.98%  "
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₃0762"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.65%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.051.55"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.37%  "
$ws.Range("E18").Value = "  +0.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "194.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.55%  "
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("E22").Value = "  -0.84%  "
$ws.Range("E23").Value = "  +5.06%  "
$ws.Range("B24").Value = "Monero"
$ws.Range("C24").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "144.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("E26").Value = "  +0.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.75%  "
$ws.Range("E28").Value = "  +0.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.24"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.28%  "
$ws.Range("E30").Value = "  -1.10%  "
$ws.Range("E31").Value = "  +1.31%  "
$ws.Range("E32").Value = "  -0.32%  "
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("E34").Value = "  +1.12%  "
$ws.Range("E35").Value = "  +0.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.132.78"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.51%  "
$ws.Range("E37").Value = "  -0.97%  "
$ws.Range("E38").Value = "  +0.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0158"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.45"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.75%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.98"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.796"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0₆0116"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "56.54"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.21%  "
$ws.Range("E45").Value = "  +3.05%  "
$ws.Range("E46").Value = "  -1.44%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.76"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.53%  "
$ws.Range("E48").Value = "  -0.14%  "
$ws.Range("E49").Value = "  +0.41%  "
$ws.Range("E50").Value = "  -1.03%  "
$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.54"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.26%  "
